# Regenerate save_data column G ("K") values - replacing old Strike# based
# values with the recalculated K counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 0
    4  = 3
    5  = 1
    6  = 8
    7  = 1
    8  = 1
    9  = 1
    10 = 4
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 2
    17 = 0
    18 = 0
    19 = 5
    20 = 2
    21 = 2
    22 = 1
    23 = 1
    24 = 0
    25 = 2
    26 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
